$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count()

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # ---- Row 8 ----
    # B8: teacher changes from "Vo Van F" to "Hoang Thi E"
    $ws.Range("B8").Value = "Hoàng Thị E"

    # D8: gains the "class" style (copy format from the stable style-8 cell H8)
    # and gets the class info that used to live in G8 (CL05 / Ky nang mem / R103)
    $ws.Range("H8").Copy()
    $ws.Range("D8").PasteSpecial(-4122)
    $ws.Range("D8").Value = "Lớp: CL05`nMôn: Kỹ năng mềm`nPhòng: R103`n(Lý thuyết)"

    # G8: loses the "class" style (copy plain format from stable style-7 cell C8)
    # and becomes empty
    $ws.Range("C8").Copy()
    $ws.Range("G8").PasteSpecial(-4122)
    $ws.Range("G8").Value = ""

    # ---- Row 9 ----
    # A9: slot time changes from C1 (13:00-15:00) to C2 (15:00-17:00)
    $ws.Range("A9").Value = "C2`n(15:00-17:00)"

    # G9: gains the "class" style and the class info that used to live in F9,
    # but now in room R103 instead of R102
    $ws.Range("H8").Copy()
    $ws.Range("G9").PasteSpecial(-4122)
    $ws.Range("G9").Value = "Lớp: CL10`nMôn: Tiếng Anh chuyên ngành`nPhòng: R103`n(Lý thuyết)"

    # F9: loses the "class" style and becomes empty
    $ws.Range("C8").Copy()
    $ws.Range("F9").PasteSpecial(-4122)
    $ws.Range("F9").Value = ""

    # ---- Row 10 ----
    # B10: teacher changes from "Ngo Van I" to "Hoang Thi E"
    $ws.Range("B10").Value = "Hoàng Thị E"

    # G10: same class/style, only the room number changes from R104 to R101
    $ws.Range("G10").Value = "Lớp: CL10`nMôn: Kỹ năng mềm`nPhòng: R101`n(Lý thuyết)"
}
